$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-07-08 Tuesday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-07-09 Wednesday", 2)

$tbl = $d.Tables.Item(1)

function Set-CellText($row, $col, $newText) {
    $cell = $tbl.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}

Set-CellText 1 1 "85÷5="
Set-CellText 1 2 "17÷8="
Set-CellText 1 3 "93÷7="
Set-CellText 1 4 "19÷6="
Set-CellText 1 5 "73÷2="

Set-CellText 5 1 "54÷6="
Set-CellText 5 2 "60÷9="
Set-CellText 5 3 "62÷3="
Set-CellText 5 4 "73÷3="
Set-CellText 5 5 "79÷6="

Set-CellText 9 1 "22÷6="
Set-CellText 9 2 "79÷4="
Set-CellText 9 3 "90÷6="
Set-CellText 9 4 "35÷5="
Set-CellText 9 5 "19÷4="

Set-CellText 13 1 "41÷8="
Set-CellText 13 2 "36÷7="
Set-CellText 13 3 "21÷6="
Set-CellText 13 4 "10÷8="
Set-CellText 13 5 "48÷8="

Set-CellText 17 1 "67÷8="
Set-CellText 17 2 "36÷3="
Set-CellText 17 3 "42÷7="
Set-CellText 17 4 "18÷6="
Set-CellText 17 5 "58÷9="
